{"js": "// Each entry is [oldText, newText]; every oldText is unique in the document,\n// so a plain (non-wildcard) case-sensitive search safely targets one run.\nconst replacements = [\n  [\"2024-10-17 Thursday\", \"2024-10-18 Friday\"],\n  [\"585\u00f78=73, 1\", \"372\u00f72=186, 0\"],\n  [\"978\u00f79=108, 6\", \"827\u00f76=137, 5\"],\n  [\"655\u00f74=163, 3\", \"759\u00f78=94, 7\"],\n  [\"532\u00f72=266, 0\", \"751\u00f74=187, 3\"],\n  [\"415\u00f76=69, 1\", \"742\u00f79=82, 4\"],\n  [\"171\u00f79=19, 0\", \"762\u00f75=152, 2\"],\n  [\"339\u00f76=56, 3\", \"203\u00f73=67, 2\"],\n  [\"794\u00f73=264, 2\", \"967\u00f76=161, 1\"],\n  [\"540\u00f73=180, 0\", \"458\u00f72=229, 0\"],\n  [\"717\u00f75=143, 2\", \"420\u00f72=210, 0\"],\n  [\"133\u00f74=33, 1\", \"560\u00f72=280, 0\"],\n  [\"852\u00f76=142, 0\", \"399\u00f79=44, 3\"],\n  [\"702\u00f73=234, 0\", \"140\u00f78=17, 4\"],\n  [\"204\u00f78=25, 4\", \"500\u00f75=100, 0\"],\n  [\"955\u00f75=191, 0\", \"706\u00f77=100, 6\"],\n  [\"837\u00f76=139, 3\", \"240\u00f75=48, 0\"],\n  [\"576\u00f79=64, 0\", \"795\u00f77=113, 4\"],\n  [\"732\u00f72=366, 0\", \"429\u00f78=53, 5\"],\n  [\"826\u00f74=206, 2\", \"156\u00f79=17, 3\"],\n  [\"436\u00f76=72, 4\", \"633\u00f75=126, 3\"],\n  [\"988\u00f74=247, 0\", \"995\u00f79=110, 5\"],\n  [\"634\u00f77=90, 4\", \"899\u00f78=112, 3\"],\n  [\"603\u00f77=86, 1\", \"393\u00f74=98, 1\"],\n  [\"238\u00f73=79, 1\", \"137\u00f74=34, 1\"],\n  [\"160\u00f73=53, 1\", \"731\u00f79=81, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Each entry is (oldText, newText); every oldText is unique in the document,\n# so Find/Replace with MatchCase + no wildcards safely targets one run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-17 Thursday\", \"2024-10-18 Friday\"),\n    @(\"585\u00f78=73, 1\", \"372\u00f72=186, 0\"),\n    @(\"978\u00f79=108, 6\", \"827\u00f76=137, 5\"),\n    @(\"655\u00f74=163, 3\", \"759\u00f78=94, 7\"),\n    @(\"532\u00f72=266, 0\", \"751\u00f74=187, 3\"),\n    @(\"415\u00f76=69, 1\", \"742\u00f79=82, 4\"),\n    @(\"171\u00f79=19, 0\", \"762\u00f75=152, 2\"),\n    @(\"339\u00f76=56, 3\", \"203\u00f73=67, 2\"),\n    @(\"794\u00f73=264, 2\", \"967\u00f76=161, 1\"),\n    @(\"540\u00f73=180, 0\", \"458\u00f72=229, 0\"),\n    @(\"717\u00f75=143, 2\", \"420\u00f72=210, 0\"),\n    @(\"133\u00f74=33, 1\", \"560\u00f72=280, 0\"),\n    @(\"852\u00f76=142, 0\", \"399\u00f79=44, 3\"),\n    @(\"702\u00f73=234, 0\", \"140\u00f78=17, 4\"),\n    @(\"204\u00f78=25, 4\", \"500\u00f75=100, 0\"),\n    @(\"955\u00f75=191, 0\", \"706\u00f77=100, 6\"),\n    @(\"837\u00f76=139, 3\", \"240\u00f75=48, 0\"),\n    @(\"576\u00f79=64, 0\", \"795\u00f77=113, 4\"),\n    @(\"732\u00f72=366, 0\", \"429\u00f78=53, 5\"),\n    @(\"826\u00f74=206, 2\", \"156\u00f79=17, 3\"),\n    @(\"436\u00f76=72, 4\", \"633\u00f75=126, 3\"),\n    @(\"988\u00f74=247, 0\", \"995\u00f79=110, 5\"),\n    @(\"634\u00f77=90, 4\", \"899\u00f78=112, 3\"),\n    @(\"603\u00f77=86, 1\", \"393\u00f74=98, 1\"),\n    @(\"238\u00f73=79, 1\", \"137\u00f74=34, 1\"),\n    @(\"160\u00f73=53, 1\", \"731\u00f79=81, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
